$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $old = $cell.Value2
    $new = $old -replace '^GUI/Examples', 'Test_Examples'
    $cell.Value2 = $new
}
